$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A11").Value = "IM2325000002"
